$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Add new worksheet right after Sheet1
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet2"

# Populate Sheet2 with the resonance data (columns A-F, rows 2-11)
$data = @(
    @(33.2, 61.1, 58.4, 16.4, 34.4, 32.6),
    @(98.4, 124.8, 182.6, 53.1, 68.7, 102),
    @(164.6, 189.3, 305.9, 88.7, 104.9, 173.1),
    @(230.7, 253.8, 430.1, 125.3, 141.2, 242.9),
    @(295.9, 319.3, $null, 162.3, 177.5, 312.6),
    @(362.1, 388, $null, 198.7, 213.2, 383.2),
    @(428.4, 457.2, $null, 235.1, 249.1, 454.2),
    @(495.9, 522.3, $null, 271.2, 285.2, 525.3),
    @($null, $null, $null, 307.6, 321.6, $null),
    @($null, $null, $null, 343.8, 358, $null)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $vals = $data[$i]
    for ($c = 0; $c -lt $vals.Length; $c++) {
        if ($null -ne $vals[$c]) {
            $ws2.Cells.Item($row, $c + 1).Value = $vals[$c]
        }
    }
}

# Set Sheet2 view: make it the selected tab with a selection of A2:F11
$ws2.Activate()
$ws2.Range("A2:F11").Select()

# Update Sheet1's view: scroll so row 3 is at top, selection B13:I13 with active cell I13
$ws1.Activate()
$ws1.Application.ActiveWindow.ScrollRow = 3
$ws1.Range("B13:I13").Select()
$excel.ActiveCell = $ws1.Range("I13")

# Re-activate Sheet2 as the final active sheet (matches tabSelected on sheet2 / activeTab=1)
$ws2.Activate()
$ws2.Range("A2:F11").Select()

# Bump the calc id slightly to mirror the saved-by-newer-Excel-build marker
$wb.ExcelApplication = $wb.ExcelApplication
